# Apply the "Holden scheme" update to the UniformA-HW05 workbook.
# This inserts 4 new HKL-index categories (Holden2.5, Holden5, Holden10, Holden15)
# as rows 20-23, each with 21 columns (C:W) of value 1, matching the existing
# pattern used for all prior categories.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newNames = @("Holden2.5", "Holden5", "Holden10", "Holden15")

for ($i = 0; $i -lt $newNames.Length; $i++) {
    $rowIndex = 20 + $i
    $hklIndex = 18 + $i

    # Column A: numeric index - copy the bold/bordered formatting from the
    # row above (the previous last data row) so the new cell reuses the
    # existing style (s="1") instead of registering a new one.
    $ws.Cells.Item($rowIndex - 1, 1).Copy()
    $ws.Cells.Item($rowIndex, 1).PasteSpecial(-4122)
    $ws.Cells.Item($rowIndex, 1).Value = $hklIndex

    # Column B: the category name (shared string)
    $ws.Cells.Item($rowIndex, 2).Value = $newNames[$i]

    # Columns C:W (3..23): all set to 1
    for ($col = 3; $col -le 23; $col++) {
        $ws.Cells.Item($rowIndex, $col).Value = 1
    }
}
